# Atualizado por script em 08-11-2023 20:45
#
# 1) Rows 16 and 17 (Chippa Utd. vs TS Galaxy / Cape Town Spurs vs Sekhukhune)
#    had their match data swapped (indices A/E stay put, F:V content swaps).
# 2) Three new match rows (79-81) are appended at the end of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Swap the F:V content of rows 16 and 17 -----------------------------

$row16 = @{
    F = "Chippa Utd."
    H = "TS Galaxy"
    I = 0
    J = 2.61
    K = "06/08/2023 21:12"
    L = 2.68
    M = "09/08/2023 19:20"
    N = 2.95
    O = "06/08/2023 21:12"
    P = 2.76
    Q = "09/08/2023 19:25"
    R = 3.1
    S = "06/08/2023 21:12"
    T = 3.24
    U = "09/08/2023 19:20"
    V = "https://www.betexplorer.com/football/south-africa/premier-league/chippa-utd-ts-galaxy/zc1AYxPj/"
}

$row17 = @{
    F = "Cape Town Spurs"
    H = "Sekhukhune"
    I = 2
    J = 3.15
    K = "07/08/2023 22:12"
    L = 2.9
    M = "09/08/2023 19:26"
    N = 2.96
    O = "07/08/2023 22:12"
    P = 2.77
    Q = "09/08/2023 19:26"
    R = 2.51
    S = "07/08/2023 22:12"
    T = 2.96
    U = "09/08/2023 19:06"
    V = "https://www.betexplorer.com/football/south-africa/premier-league/cape-town-spurs-sekhukhune/SjlddtHc/"
}

foreach ($col in $row17.Keys) {
    $ws.Range($col + "16").Value = $row17[$col]
}
foreach ($col in $row16.Keys) {
    $ws.Range($col + "17").Value = $row16[$col]
}

# --- 2) Append new rows 79, 80, 81 ------------------------------------------

$newRows = @(
    @{
        Row = 79
        A = 78
        B = "south-africa"
        C = "premier-league"
        D = "2023-2024"
        E = 45238.77083333334
        F = "Kaizer Chiefs"
        G = 3
        H = "Cape Town Spurs"
        I = 2
        J = 1.62
        K = "01/11/2023 18:43"
        L = 1.44
        M = "08/11/2023 18:20"
        N = 3.63
        O = "01/11/2023 18:43"
        P = 4.21
        Q = "08/11/2023 18:29"
        R = 6.05
        S = "01/11/2023 18:43"
        T = 8.36
        U = "08/11/2023 18:25"
        V = "https://www.betexplorer.com/football/south-africa/premier-league/kaizer-chiefs-cape-town-spurs/8O01wsT8/"
    },
    @{
        Row = 80
        A = 79
        B = "south-africa"
        C = "premier-league"
        D = "2023-2024"
        E = 45238.77083333334
        F = "Polokwane"
        G = 1
        H = "Supersport Utd"
        I = 1
        J = 2.77
        K = "06/11/2023 12:42"
        L = 3.73
        M = "08/11/2023 18:26"
        N = 3.01
        O = "06/11/2023 12:42"
        P = 2.83
        Q = "08/11/2023 18:26"
        R = 2.65
        S = "06/11/2023 12:42"
        T = 2.36
        U = "08/11/2023 18:26"
        V = "https://www.betexplorer.com/football/south-africa/premier-league/polokwane-city-supersport-utd/2o05x1rF/"
    },
    @{
        Row = 81
        A = 80
        B = "south-africa"
        C = "premier-league"
        D = "2023-2024"
        E = 45238.77083333334
        F = "Swallows"
        G = 1
        H = "TS Galaxy"
        I = 0
        J = 2.56
        K = "06/11/2023 12:42"
        L = 2.41
        M = "08/11/2023 18:13"
        N = 2.83
        O = "06/11/2023 12:42"
        P = 2.81
        Q = "08/11/2023 18:13"
        R = 3.22
        S = "06/11/2023 12:42"
        T = 3.67
        U = "08/11/2023 18:13"
        V = "https://www.betexplorer.com/football/south-africa/premier-league/swallows-fc-ts-galaxy/txa9yLcL/"
    }
)

$lastRow = 78
$colOrder = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

foreach ($rowData in $newRows) {
    $r = $rowData.Row

    # Copy formatting (styles) from the previous row first.
    $ws.Range("A" + $lastRow + ":V" + $lastRow).Copy()
    $ws.Range("A" + $r + ":V" + $r).PasteSpecial(-4122)

    foreach ($col in $colOrder) {
        $ws.Range($col + $r).Value = $rowData[$col]
    }

    $lastRow = $r
}

$excel.CutCopyMode = 0
